$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3150
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -4350
$ws.Range("H43").Value = 6950919.5
$ws.Range("I43").Value = 20300.5
$ws.Range("J43").Value = 9261126
$ws.Range("K43").Value = 20300.5
$ws.Range("L43").Value = 9261126
$ws.Range("M43").Value = -20231.5
$ws.Range("N43").Value = -9261264
$ws.Range("H111").Value = 4102.75
$ws.Range("I111").Value = 4747
$ws.Range("J111").Value = 3029
$ws.Range("K111").Value = 14241
$ws.Range("L111").Value = 9087
$ws.Range("M111").Value = -11174
$ws.Range("N111").Value = -15221
$ws.Range("H129").Value = 828.23883
$ws.Range("I129").Value = 395
$ws.Range("J129").Value = 913.3393
$ws.Range("K129").Value = 1185
$ws.Range("L129").Value = 2740.0179
$ws.Range("M129").Value = 3815
$ws.Range("N129").Value = -12740.0179
$ws.Range("H135").Value = 279.75
$ws.Range("I135").Value = 235.7
$ws.Range("K135").Value = 2121.3
$ws.Range("M135").Value = 413.7000000000003
$ws.Range("H137").Value = 1049.9688
$ws.Range("J137").Value = 1250.6
$ws.Range("L137").Value = 3751.8
$ws.Range("N137").Value = -8851.799999999999
$ws.Range("H138").Value = 1353.5278
$ws.Range("I138").Value = 728.9655
$ws.Range("J138").Value = 1774.7441
$ws.Range("K138").Value = 2186.8965
$ws.Range("L138").Value = 5324.2323
$ws.Range("M138").Value = 2953.1035
$ws.Range("N138").Value = -15604.2323
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H11").Value = 5449.75
$ws.Range("I11").Value = 5599.6665
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 5599.6665
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -5455.6665
$ws.Range("N11").Value = -5288
$ws.Range("H32").Value = 4950.2856
$ws.Range("I32").Value = 5526.1665
$ws.Range("J32").Value = 1495
$ws.Range("K32").Value = 5526.1665
$ws.Range("L32").Value = 1495
$ws.Range("M32").Value = -5239.1665
$ws.Range("N32").Value = -2069
$ws.Range("H97").Value = 511.25
$ws.Range("I97").Value = 381.66666
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 381.66666
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = 114.33334
$ws.Range("N97").Value = -1892
$ws.Range("H102").Value = 27778596
$ws.Range("I102").Value = 27778596
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 27778596
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -27776974
$ws.Range("N102").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 25001038
$ws.Range("I99").Value = 31250906
$ws.Range("J99").Value = 1560.25
$ws.Range("K99").Value = 31250906
$ws.Range("L99").Value = 1560.25
$ws.Range("M99").Value = -31249408
$ws.Range("N99").Value = -4556.25
$ws.Range("H134").Value = 6275.905
$ws.Range("I134").Value = 1549.6875
$ws.Range("K134").Value = 4649.0625
$ws.Range("M134").Value = -2114.0625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8064.9414
$ws.Range("I132").Value = 12223.444
$ws.Range("K132").Value = 36670.33199999999
$ws.Range("M132").Value = -34140.33199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1597
$ws.Range("I34").Value = 673.3333
$ws.Range("J34").Value = 1992.8572
$ws.Range("K34").Value = 2019.9999
$ws.Range("L34").Value = 5978.571599999999
$ws.Range("M34").Value = -1935.9999
$ws.Range("N34").Value = -6146.571599999999
$ws.Range("H139").Value = 3929
$ws.Range("I139").Value = 4502.25
$ws.Range("J139").Value = 3470.4
$ws.Range("K139").Value = 13506.75
$ws.Range("L139").Value = 10411.2
$ws.Range("M139").Value = -8366.75
$ws.Range("N139").Value = -20691.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6584375
$ws.Range("I11").Value = 6944643
$ws.Range("K11").Value = 6944643
$ws.Range("M11").Value = -6944504
$ws.Range("H70").Value = 20460632
$ws.Range("I70").Value = 19236094
$ws.Range("J70").Value = 22229412
$ws.Range("K70").Value = 19236094
$ws.Range("L70").Value = 22229412
$ws.Range("M70").Value = -19235824
$ws.Range("N70").Value = -22229952
$ws.Range("H73").Value = 20460632
$ws.Range("I73").Value = 19236094
$ws.Range("J73").Value = 22229412
$ws.Range("K73").Value = 19236094
$ws.Range("L73").Value = 22229412
$ws.Range("M73").Value = -19235158
$ws.Range("N73").Value = -22231284
$ws.Range("H132").Value = 2513.2
$ws.Range("I132").Value = 2191.1333
$ws.Range("K132").Value = 6573.3999
$ws.Range("M132").Value = -4043.3999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 113
$ws.Range("I32").Value = 113
$ws.Range("K32").Value = 113
$ws.Range("M32").Value = 204
$ws.Range("H61").Value = 1291.6
$ws.Range("J61").Value = 901.6667
$ws.Range("L61").Value = 901.6667
$ws.Range("N61").Value = -1305.6667
$ws.Range("H113").Value = 1291.6
$ws.Range("J113").Value = 901.6667
$ws.Range("L113").Value = 901.6667
$ws.Range("N113").Value = -5241.6667
$ws.Range("H132").Value = 79907.766
$ws.Range("I132").Value = 3699.5
$ws.Range("K132").Value = 11098.5
$ws.Range("M132").Value = -8568.5
$ws.Range("H136").Value = 7139.722
$ws.Range("I136").Value = 8939.691999999999
$ws.Range("K136").Value = 26819.076
$ws.Range("M136").Value = -24269.076
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2333.3333
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2333.3333
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2333.3333
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = -2559.3333
$ws.Range("H122").Value = 11820550
$ws.Range("I122").Value = 13002510
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 39007530
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -39005080
$ws.Range("N122").Value = -7750
$ws.Range("H136").Value = 652.5
$ws.Range("I136").Value = 460
$ws.Range("J136").Value = 1076
$ws.Range("K136").Value = 1380
$ws.Range("L136").Value = 3228
$ws.Range("M136").Value = 1170
$ws.Range("N136").Value = -8328
